$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert four new "Duy" hour-tracking rows into the log, shifting the
# existing rows (and the trailing Total row/formula) downward. Inserted
# from top to bottom so each target row number already accounts for the
# shift caused by the previous insert.

# New row 6: Duy, 2022-02-25 (44617), 4 hours
$ws.Rows("6:6").Insert()
$ws.Range("A6").Value = "Duy"
$ws.Range("B6").Value = 44617
$ws.Range("C6").Value = 4

# New row 10: Duy, 2022-03-02 (44622), 3 hours
$ws.Rows("10:10").Insert()
$ws.Range("A10").Value = "Duy"
$ws.Range("B10").Value = 44622
$ws.Range("C10").Value = 3

# New row 13: Duy, 2022-03-03 (44623), 3 hours
$ws.Rows("13:13").Insert()
$ws.Range("A13").Value = "Duy"
$ws.Range("B13").Value = 44623
$ws.Range("C13").Value = 3

# New row 15: Duy, 2022-03-04 (44624), 2 hours
$ws.Rows("15:15").Insert()
$ws.Range("A15").Value = "Duy"
$ws.Range("B15").Value = 44624
$ws.Range("C15").Value = 2

# Update the cursor/selection to where the last edit was made
[void]$ws.Range("D21").Select()
